$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (06-nov) before the
#     existing "01-oct." column block (DI), shifting DI:EM -> DJ:EN.
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("DI1").EntireColumn.Insert()

$ws1.Range("DI1").Value = "06-nov"
for ($r = 2; $r -le 25; $r++) {
  $ws1.Cells.Item($r, 113).Value = "-"
}

# --- Sheet "Gaz": append new row with latest price date.
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A142").NumberFormat = "@"
$ws2.Range("A142").Value = "2025-11-04"
$ws2.Range("B142").Value = 31.17

# --- Sheet "CO2": append new row with latest price date.
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A142").NumberFormat = "@"
$ws3.Range("A142").Value = "2025-11-04"
$ws3.Range("B142").Value = 81.9
